$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16's "Test R2" cell (F16) was an empty placeholder cell in the source
# data; the refreshed export no longer emits it at all, so clear it out
# entirely (not just blank its displayed text).
$ws.Range("F16").ClearContents()

# Append the newly logged evaluation run as row 17.
$ws.Range("A17").Value = "2024-11-26 15:17:21"
$ws.Range("B17").Value = 0.9961532047459436
$ws.Range("C17").Value = 0.00753433642676356
$ws.Range("D17").Value = 0.0001750359918135414
$ws.Range("E17").Value = 0.01323011684806833
# "Test R2" (F17) is blank for this run too - write it as an actual empty
# text cell (matching the existing blank-F-cell convention used elsewhere
# in this sheet) rather than leaving the cell completely absent.
$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = 0.02188858299558192
$ws.Range("H17").Value = 0.0005584442789189223
$ws.Range("I17").Value = 0.02363142566412197
